# Apply updated dSF (column F) values to the relevant rows in Sheet1.
# These correspond to a "repull data, push all data, mean calculation" update
# where only the F column (dSF) values were recalculated/replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -4
    4  = 4
    6  = 5
    11 = -2
    14 = -3
    19 = -4
    20 = -5
    21 = -3
    22 = -5
    27 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
